$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G to hold the new "Major Exporter's
# Production" variable, pushing the existing Global Production column
# (and the trailing note column) one slot to the right.
$ws.Columns("G:G").Insert()

# Populate the header for the newly inserted column.
$ws.Range("G1").Value = "% Change in Major Exporter's Production"

# Resize the new column to fit its (longer) header text, matching the
# width Excel's own best-fit would compute for this text.
$ws.Columns("G:G").ColumnWidth = 34.25

# Reflect the post-edit cell selection/view state.
$ws.Range("I4").Select()
